$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: insert two new columns before D, shifting existing D:K data to F:M
$ws.Range("D:E").EntireColumn.Insert()

# Step 2: copy number formats from the (now shifted) F:G columns into new D:E columns
# (only for the row blocks that actually carry quarterly data)
$ws.Range("F7:G35").Copy()
$ws.Range("D7").PasteSpecial(-4122)
$ws.Range("F38:G77").Copy()
$ws.Range("D38").PasteSpecial(-4122)
$ws.Range("F80:G102").Copy()
$ws.Range("D80").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Step 3: populate the two new quarter columns (D, E) with their reported values
$ws.Range("D7").Value = 43465
$ws.Range("E7").Value = 43373
$ws.Range("D8").Value = 93000
$ws.Range("E8").Value = 77700
$ws.Range("D9").Value = 69200
$ws.Range("E9").Value = 57800
$ws.Range("D10").Value = 23800
$ws.Range("E10").Value = 19900
$ws.Range("D12").Value = "NA"
$ws.Range("E12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("D14").Value = "NA"
$ws.Range("E14").Value = "NA"
$ws.Range("D15").Value = "NA"
$ws.Range("E15").Value = "NA"
$ws.Range("D17").Value = 88400
$ws.Range("E17").Value = 77000
$ws.Range("D18").Value = 4600
$ws.Range("E18").Value = 700
$ws.Range("D20").Value = -700
$ws.Range("E20").Value = -600
$ws.Range("D21").Value = 6600
$ws.Range("E21").Value = 2800
$ws.Range("D22").Value = 0
$ws.Range("E22").Value = 0
$ws.Range("D23").Value = 3900
$ws.Range("E23").Value = 100
$ws.Range("D24").Value = 500
$ws.Range("E24").Value = 800
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("D26").Value = 3500
$ws.Range("E26").Value = -700
$ws.Range("D27").Value = 3500
$ws.Range("E27").Value = -700
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0
$ws.Range("D29").Value = -200
$ws.Range("E29").Value = "NA"
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0
$ws.Range("D32").Value = 700
$ws.Range("E32").Value = 600
$ws.Range("D33").Value = 3300
$ws.Range("E33").Value = -700
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0
$ws.Range("D35").Value = 3300
$ws.Range("E35").Value = -700
$ws.Range("D38").Value = 43465
$ws.Range("E38").Value = 43373
$ws.Range("D41").Value = 15400
$ws.Range("E41").Value = 14300
$ws.Range("D42").Value = 0
$ws.Range("E42").Value = 0
$ws.Range("D43").Value = 72400
$ws.Range("E43").Value = 56300
$ws.Range("D44").Value = 10100
$ws.Range("E44").Value = 5200
$ws.Range("D45").Value = 7200
$ws.Range("E45").Value = 5200
$ws.Range("D46").Value = 105100
$ws.Range("E46").Value = 81000
$ws.Range("D47").Value = 0
$ws.Range("E47").Value = 0
$ws.Range("D48").Value = 21500
$ws.Range("E48").Value = 21900
$ws.Range("D49").Value = 47000
$ws.Range("E49").Value = 47500
$ws.Range("D50").Value = 0
$ws.Range("E50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 0
$ws.Range("D52").Value = 3500
$ws.Range("E52").Value = 3500
$ws.Range("D53").Value = 0
$ws.Range("E53").Value = 0
$ws.Range("D54").Value = 177100
$ws.Range("E54").Value = 153900
$ws.Range("D57").Value = 47600
$ws.Range("E57").Value = 29300
$ws.Range("D58").Value = 2600
$ws.Range("E58").Value = 3100
$ws.Range("D59").Value = 32000
$ws.Range("E59").Value = 29500
$ws.Range("D60").Value = 82100
$ws.Range("E60").Value = 61800
$ws.Range("D61").Value = 39300
$ws.Range("E61").Value = 40100
$ws.Range("D62").Value = 9000
$ws.Range("E62").Value = 9700
$ws.Range("D63").Value = 0
$ws.Range("E63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("E64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("E65").Value = 0
$ws.Range("D66").Value = 130500
$ws.Range("E66").Value = 111700
$ws.Range("D68").Value = 0
$ws.Range("E68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("E69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("E70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("E71").Value = 0
$ws.Range("D72").Value = -107800
$ws.Range("E72").Value = -111100
$ws.Range("D73").Value = 0
$ws.Range("E73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("E74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("E75").Value = 0
$ws.Range("D76").Value = 46600
$ws.Range("E76").Value = 42200
$ws.Range("D77").Value = 0
$ws.Range("E77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("E80").Value = 43373
$ws.Range("D81").Value = 3300
$ws.Range("E81").Value = -700
$ws.Range("D83").Value = 2700
$ws.Range("E83").Value = 2700
$ws.Range("D84").Value = 0
$ws.Range("E84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("E85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("E86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("E87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("E88").Value = 0
$ws.Range("D89").Value = 1100
$ws.Range("E89").Value = 7800
$ws.Range("D91").Value = -900
$ws.Range("E91").Value = -2100
$ws.Range("D92").Value = 0
$ws.Range("E92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("E93").Value = 0
$ws.Range("D94").Value = -900
$ws.Range("E94").Value = -2100
$ws.Range("D96").Value = 0
$ws.Range("E96").Value = 0
$ws.Range("D97").Value = 0
$ws.Range("E97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("E98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("E99").Value = 0
$ws.Range("D100").Value = 600
$ws.Range("E100").Value = -5000
$ws.Range("D101").Value = 400
$ws.Range("E101").Value = -100
$ws.Range("D102").Value = 1200
$ws.Range("E102").Value = 600
